$wb = $excel.ActiveWorkbook

# Updated "想去人数" (F column) values for sheets "展览" (row-level event listing)
# and "全部类型" (same data replicated). Row -> new value mapping.
$updates = @{
    2  = 1078
    3  = 783
    4  = 267
    5  = 39
    8  = 1917
    9  = 7103
    10 = 493
    11 = 393
    12 = 327
    13 = 115
    15 = 147
    16 = 7036
    18 = 1317
    19 = 142
    21 = 226
    22 = 122
    23 = 289
    24 = 123
    26 = 13
    28 = 16
    35 = 34
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
